$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; existing rows 6..45 shift down to 7..46
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with a new weekly record
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44490
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 300000001
$ws.Cells.Item(6, 7).Value = "Rabanito"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 65
$ws.Cells.Item(6, 11).Value = 6000
$ws.Cells.Item(6, 12).Value = 6000
$ws.Cells.Item(6, 13).Value = 6000
$ws.Cells.Item(6, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(6, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(6, 16).Value = 500
$ws.Cells.Item(6, 17).Value = 12
$ws.Cells.Item(6, 18).Value = "Hortaliza"
